$wb = $excel.ActiveWorkbook

# --- Sheet1: summary table (rows 2-11) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B2").Value = -152.6915513587452
$ws1.Range("C2").Value = 9.981462678
$ws1.Range("F2").Value = 20
$ws1.Range("G2").Value = 4540
$ws1.Range("H2").Value = 4900
$ws1.Range("I2").Value = 400
$ws1.Range("B3").Value = -149.80697334593236
$ws1.Range("C3").Value = 1.670559779
$ws1.Range("F3").Value = 20
$ws1.Range("G3").Value = 4540
$ws1.Range("H3").Value = 4900
$ws1.Range("I3").Value = 400
$ws1.Range("B4").Value = -151.10442773297382
$ws1.Range("C4").Value = 2.125055329
$ws1.Range("F4").Value = 20
$ws1.Range("G4").Value = 4540
$ws1.Range("H4").Value = 4900
$ws1.Range("I4").Value = 400
$ws1.Range("B5").Value = -150.7026611363612
$ws1.Range("C5").Value = 1.859133424
$ws1.Range("F5").Value = 20
$ws1.Range("G5").Value = 4540
$ws1.Range("H5").Value = 4900
$ws1.Range("I5").Value = 400
$ws1.Range("B6").Value = -150.52128244649396
$ws1.Range("C6").Value = 2.881476532
$ws1.Range("F6").Value = 20
$ws1.Range("G6").Value = 4540
$ws1.Range("H6").Value = 4900
$ws1.Range("I6").Value = 400
$ws1.Range("B7").Value = -149.30562494657448
$ws1.Range("C7").Value = 2.936565255
$ws1.Range("F7").Value = 20
$ws1.Range("G7").Value = 4540
$ws1.Range("H7").Value = 4900
$ws1.Range("I7").Value = 400
$ws1.Range("B8").Value = -148.66151974917545
$ws1.Range("C8").Value = 1.676146803
$ws1.Range("F8").Value = 20
$ws1.Range("G8").Value = 4540
$ws1.Range("H8").Value = 4900
$ws1.Range("I8").Value = 400
$ws1.Range("B9").Value = -150.18259972871869
$ws1.Range("C9").Value = 1.963609485
$ws1.Range("F9").Value = 20
$ws1.Range("G9").Value = 4540
$ws1.Range("H9").Value = 4900
$ws1.Range("I9").Value = 400
$ws1.Range("B10").Value = -151.08461007105797
$ws1.Range("C10").Value = 2.07834734
$ws1.Range("F10").Value = 20
$ws1.Range("G10").Value = 4540
$ws1.Range("H10").Value = 4900
$ws1.Range("I10").Value = 400
$ws1.Range("B11").Value = -146.922210513398
$ws1.Range("C11").Value = 1.777984944
$ws1.Range("F11").Value = 20
$ws1.Range("G11").Value = 4540
$ws1.Range("H11").Value = 4900
$ws1.Range("I11").Value = 400

# --- Tab "1": CCG iteration log ---
$wsT = $wb.Worksheets.Item("1")
$wsT.Range("D2").Value = 0.885273730232666
$wsT.Range("E2").Value = 220.70464
$wsT.Range("B3").Value = -152.6915513587452
$wsT.Range("C3").Value = 0.09094437275100899
$wsT.Range("D3").Value = 1.3182471356051026

# --- Tab "2": CCG iteration log ---
$wsT = $wb.Worksheets.Item("2")
$wsT.Range("D2").Value = 0.019672506603149415
$wsT.Range("E2").Value = 222.6832
$wsT.Range("B3").Value = -149.80697334593236
$wsT.Range("C3").Value = 0.04863840334435918
$wsT.Range("D3").Value = 1.4584250013929443

# --- Tab "3": CCG iteration log ---
$wsT = $wb.Worksheets.Item("3")
$wsT.Range("D2").Value = 0.011684348206420898
$wsT.Range("E2").Value = 222.06669
$wsT.Range("B3").Value = -151.10442773297382
$wsT.Range("D3").Value = 1.9546205120506592

# --- Tab "4": CCG iteration log ---
$wsT = $wb.Worksheets.Item("4")
$wsT.Range("D2").Value = 0.03785099905883789
$wsT.Range("E2").Value = 219.63896
$wsT.Range("B3").Value = -150.7026611363612
$wsT.Range("C3").Value = 0.0790679019001801
$wsT.Range("D3").Value = 1.6570699790717773

# --- Tab "5": CCG iteration log ---
$wsT = $wb.Worksheets.Item("5")
$wsT.Range("D2").Value = 0.0369043198894043
$wsT.Range("E2").Value = 226.07773
$wsT.Range("B3").Value = -150.52128244649396
$wsT.Range("D3").Value = 2.632920491957031

# --- Tab "6": CCG iteration log ---
$wsT = $wb.Worksheets.Item("6")
$wsT.Range("D2").Value = 0.046772064153686524
$wsT.Range("E2").Value = 236.63439
$wsT.Range("B3").Value = -149.30562494657448
$wsT.Range("C3").Value = 0.0
$wsT.Range("D3").Value = 2.5871455994648436

# --- Tab "7": CCG iteration log ---
$wsT = $wb.Worksheets.Item("7")
$wsT.Range("D2").Value = 0.03861509391125488
$wsT.Range("E2").Value = 217.43514
$wsT.Range("B3").Value = -148.66151974917545
$wsT.Range("D3").Value = 1.416365807449829

# --- Tab "8": CCG iteration log ---
$wsT = $wb.Worksheets.Item("8")
$wsT.Range("D2").Value = 0.013593793322265625
$wsT.Range("E2").Value = 213.70945
$wsT.Range("B3").Value = -150.18259972871869
$wsT.Range("C3").Value = 0.017091668477038063
$wsT.Range("D3").Value = 1.7661272293115235

# --- Tab "9": CCG iteration log ---
$wsT = $wb.Worksheets.Item("9")
$wsT.Range("D2").Value = 0.028167553452148438
$wsT.Range("E2").Value = 228.34639
$wsT.Range("B3").Value = -151.08461007105797
$wsT.Range("C3").Value = 0.05281059687761344
$wsT.Range("D3").Value = 1.8342621356691895

# --- Tab "10": CCG iteration log ---
$wsT = $wb.Worksheets.Item("10")
$wsT.Range("D2").Value = 0.03731891070141601
$wsT.Range("E2").Value = 218.51635
$wsT.Range("B3").Value = -146.922210513398
$wsT.Range("C3").Value = 0.008927839974904692
$wsT.Range("D3").Value = 1.5419986081529542

